# Apply the data refresh for Actual_Production_Solar.xlsx
# 1) Shift every timestamp in column A (rows 2-97) forward by exactly one day.
# 2) Update the solar production values in column B for rows 24-45 with the
#    newer model's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift all timestamps in column A by +1 day ---
$lastRow = 97
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# --- Step 2: update column B values for rows 24-45 ---
$newB = @{
    24 = 13
    25 = 34
    26 = 79
    27 = 129
    28 = 173
    29 = 239
    30 = 341
    31 = 431
    32 = 496
    33 = 486
    34 = 581
    35 = 657
    36 = 671
    37 = 723
    38 = 790
    39 = 797
    40 = 850
    41 = 868
    42 = 880
    43 = 930
    44 = 945
    45 = 944
}

foreach ($row in $newB.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newB[$row]
}
